$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns D and E stay as text (avoid Excel auto-converting
# numeric-looking strings like "1.00" or "152.40" into numbers).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.659.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.24%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.613.95'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.96%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.53'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.19%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.40'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.39%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.543'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.24%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.612.75'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.03%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.76%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.72%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.19'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.60%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.346'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.80%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.50'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.94%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000189'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.52%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.092.22'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.95%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.554.37'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.05%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.612.12'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '371.81'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.92%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.18'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.20%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.05'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -7.71%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.00%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.31%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.25%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.83'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +10.12%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.02%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.86'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.14%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '596.08'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.25%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.748.29'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.78%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000104'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.95%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.04%  '

# Row 32
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.38'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.08%  '

# Row 33
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.79'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.71%  '

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.03%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.125'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.95%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.51'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.03%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.43'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.56%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.13'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.85%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.78%  '

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.44%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.26'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.08%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.71'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.49%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.10'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.51%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.08%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.39'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.12%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '156.05'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.35%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0296'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.35%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.68%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.69'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.05%  '

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.73%  '
